# Generate Report for handoff
#
# Updates the "Latest Handoff Datetime" (column D) for the e398423b...
# row (row 5) on both the "zh-cn" and "de-de" localization-status sheets,
# recording the timestamp at which this handoff report was generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-19 05:07:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-19 05:07:44"
